$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Upgrade utility poles from wood to metal" row (row 2).
# This shifts the remaining option rows up by one.
$ws.Rows(2).Delete()

# Add a new option row for substation flood protection (row 9, leaving a gap at 7-8).
$ws.Range("A9").Value = 'Substation flood protection'
$ws.Range("B9").Value = 'junction'
$ws.Range("C9").Value = 'substation'
$ws.Range("D9").Value = 4289168
$ws.Range("E9").Value = '$USD per substation'
$ws.Range("F9").Value = 'Full protection from flooding'
$ws.Range("G9").Value = 'Based on National Grid analysis from the UK'
$ws.Hyperlinks.Add($ws.Range("H9"), "https://www.nationalgrid.com/uk/electricity-transmission/document/140911/download")

# Re-add the utility pole upgrade option (now "wood to steel") with updated cost data (row 10).
$ws.Range("A10").Value = 'Upgrade utility poles from wood to steel '
$ws.Range("B10").Value = 'junction'
$ws.Range("C10").Value = 'pole'
$ws.Range("D10").Value = 800.19
$ws.Range("E10").Value = '$USD/kV per asset'
$ws.Range("F10").Value = 'Results in a 76% reduction in damage index (Miyamoto 2019)'
$ws.Range("G10").Value = 'Cost data based on MISO Energy (U.S.)'
$ws.Range("H10").Value = 'https://cdn.misoenergy.org/20190212%20PSC%20Item%2005a%20Transmission%20Cost%20Estimation%20Guide%20for%20MTEP%202019_for%20review317692.pdf'

# Column width adjustments for the new columns.
$ws.Columns("D").ColumnWidth = 8.6640625
$ws.Columns("E").ColumnWidth = 18
$ws.Columns("H").ColumnWidth = 75.6640625

# View adjustments.
$ws.Range("A2:XFD2").Select()
$excel.ActiveWindow.Zoom = 132
